$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 563.3
$ws.Range("I2").Value = 591.625
$ws.Range("K2").Value = 591.625
$ws.Range("M2").Value = -478.625

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 177.45454
$ws.Range("I33").Value = 187.875
$ws.Range("J33").Value = 149.66667
$ws.Range("K33").Value = 187.875
$ws.Range("L33").Value = 149.66667
$ws.Range("M33").Value = 41.125
$ws.Range("N33").Value = -607.6666700000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 66.25
$ws.Range("J42").Value = 100
$ws.Range("L42").Value = 300
$ws.Range("N42").Value = -760

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1749.3334
$ws.Range("I70").Value = 499
$ws.Range("K70").Value = 1497
$ws.Range("M70").Value = -1227

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1749.3334
$ws.Range("I73").Value = 499
$ws.Range("K73").Value = 1497
$ws.Range("M73").Value = -561

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 9000
$ws.Range("J116").Value = 9000
$ws.Range("L116").Value = 9000
$ws.Range("N116").Value = -15884

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5250.25
$ws.Range("I137").Value = 6666
$ws.Range("J137").Value = 1003
$ws.Range("K137").Value = 19998
$ws.Range("L137").Value = 3009
$ws.Range("M137").Value = -17448
$ws.Range("N137").Value = -8109

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 2820
$ws.Range("J14").Value = 700
$ws.Range("L14").Value = 700
$ws.Range("N14").Value = -1050

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 5000
$ws.Range("J19").Value = 5000
$ws.Range("L19").Value = 5000
$ws.Range("N19").Value = -5458

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 100
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 100
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -848

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6398
$ws.Range("I61").Value = 1330
$ws.Range("K61").Value = 1330
$ws.Range("M61").Value = -1118

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9750
$ws.Range("I74").Value = 3000
$ws.Range("K74").Value = 3000
$ws.Range("M74").Value = -2126

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 9750
$ws.Range("I77").Value = 3000
$ws.Range("K77").Value = 15000
$ws.Range("M77").Value = -10632

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1000
$ws.Range("I110").Value = 1000
$ws.Range("K110").Value = 1000
$ws.Range("M110").Value = 1045

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3404.5625
$ws.Range("I132").Value = 1047.75
$ws.Range("J132").Value = 10475
$ws.Range("K132").Value = 3143.25
$ws.Range("L132").Value = 31425
$ws.Range("M132").Value = -613.25
$ws.Range("N132").Value = -36485

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 6398
$ws.Range("I136").Value = 1330
$ws.Range("K136").Value = 3990
$ws.Range("M136").Value = -1440

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 6500
$ws.Range("J23").Value = 6500
$ws.Range("L23").Value = 6500
$ws.Range("N23").Value = -7066

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 5000
$ws.Range("J30").Value = 5000
$ws.Range("L30").Value = 5000
$ws.Range("N30").Value = -5250

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 5525.8335
$ws.Range("I82").Value = 5525.8335
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 5525.8335
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -5142.8335
$ws.Range("N82").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H85").Value = 5525.8335
$ws.Range("I85").Value = 5525.8335
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 5525.8335
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -4199.8335
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 28380
$ws.Range("I102").Value = 16725
$ws.Range("K102").Value = 16725
$ws.Range("M102").Value = -13480

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2086.3
$ws.Range("I105").Value = 2040.3334
$ws.Range("K105").Value = 2040.3334
$ws.Range("M105").Value = -293.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5555.5557
$ws.Range("I4").Value = 5000
$ws.Range("K4").Value = 5000
$ws.Range("M4").Value = -4888

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5854.1177
$ws.Range("I31").Value = 2012.1
$ws.Range("J31").Value = 11342.714
$ws.Range("K31").Value = 2012.1
$ws.Range("L31").Value = 11342.714
$ws.Range("M31").Value = -1717.1
$ws.Range("N31").Value = -11932.714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5854.1177
$ws.Range("I34").Value = 2012.1
$ws.Range("J34").Value = 11342.714
$ws.Range("K34").Value = 2012.1
$ws.Range("L34").Value = 11342.714
$ws.Range("M34").Value = -1810.1
$ws.Range("N34").Value = -11746.714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 6832.231
$ws.Range("I58").Value = 3868.889
$ws.Range("K58").Value = 3868.889
$ws.Range("M58").Value = -3665.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 11928.571
$ws.Range("J86").Value = 12500
$ws.Range("L86").Value = 12500
$ws.Range("N86").Value = -14746

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 11928.571
$ws.Range("J89").Value = 12500
$ws.Range("L89").Value = 62500
$ws.Range("N89").Value = -73732

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 8819.799999999999
$ws.Range("I132").Value = 3039.6
$ws.Range("K132").Value = 9118.799999999999
$ws.Range("M132").Value = -6588.799999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4303.75
$ws.Range("I134").Value = 1518.7142
$ws.Range("K134").Value = 4556.142599999999
$ws.Range("M134").Value = -2021.142599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 6832.231
$ws.Range("I136").Value = 3868.889
$ws.Range("K136").Value = 11606.667
$ws.Range("M136").Value = -9056.667000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 545815
$ws.Range("J141").Value = 545815
$ws.Range("L141").Value = 545815
$ws.Range("N141").Value = -556175

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 129
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 3277.5
$ws.Range("I122").Value = 1555
$ws.Range("K122").Value = 13995
$ws.Range("M122").Value = -11545

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1000
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -1224

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4623.643
$ws.Range("I132").Value = 2603.1
$ws.Range("J132").Value = 9675
$ws.Range("K132").Value = 7809.299999999999
$ws.Range("L132").Value = 29025
$ws.Range("M132").Value = -5279.299999999999
$ws.Range("N132").Value = -34085

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1686.375
$ws.Range("I61").Value = 1398.8
$ws.Range("K61").Value = 1398.8
$ws.Range("M61").Value = -1196.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1686.375
$ws.Range("I113").Value = 1398.8
$ws.Range("K113").Value = 1398.8
$ws.Range("M113").Value = 771.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3555.7144
$ws.Range("I122").Value = 3231.6667
$ws.Range("J122").Value = 5500
$ws.Range("K122").Value = 9695.000100000001
$ws.Range("L122").Value = 16500
$ws.Range("M122").Value = -7245.000100000001
$ws.Range("N122").Value = -21400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 7450.8
$ws.Range("I132").Value = 5396.4165
$ws.Range("K132").Value = 16189.2495
$ws.Range("M132").Value = -13659.2495

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 11426.857
$ws.Range("J136").Value = 13998
$ws.Range("L136").Value = 41994
$ws.Range("N136").Value = -47094

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 784.6667
$ws.Range("I11").Value = 227
$ws.Range("J11").Value = 1900
$ws.Range("K11").Value = 227
$ws.Range("L11").Value = 1900
$ws.Range("M11").Value = -85
$ws.Range("N11").Value = -2184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 6749.5
$ws.Range("I34").Value = 6749.5
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 6749.5
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -6546.5
$ws.Range("N34").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 81666
$ws.Range("J46").Value = 81666
$ws.Range("L46").Value = 81666
$ws.Range("N46").Value = -82128

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 200000
$ws.Range("J62").Value = 200000
$ws.Range("L62").Value = 200000
$ws.Range("N62").Value = -201248

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 200000
$ws.Range("J65").Value = 200000
$ws.Range("L65").Value = 1000000
$ws.Range("N65").Value = -1006240

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H134").Value = 81666
$ws.Range("J134").Value = 81666
$ws.Range("L134").Value = 244998
$ws.Range("N134").Value = -250068

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 8215.846
$ws.Range("I136").Value = 6800.5454
$ws.Range("K136").Value = 20401.6362
$ws.Range("M136").Value = -17851.6362
